$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Locate the trailing paragraphs at the end of the document body:
#    ... "{%p endfor %}"  (target - keep, but reformat)
#        "changeed"       (delete)
#        ""               (delete - trailing empty paragraph)
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$changeedPara = $lastPara.Previous()
$targetPara = $changeedPara.Previous()

# sanity checks (no-ops if already correct, just documents intent)
# $targetPara.Range.Text   -> "{%p endfor %}"
# $changeedPara.Range.Text -> "changeed"
# $lastPara.Range.Text     -> ""

# ------------------------------------------------------------------
# 2. Remove the "changeed" paragraph and the empty paragraph after it
#    (delete from the start of the "changeed" paragraph through the
#    end of the final empty paragraph).
# ------------------------------------------------------------------
$delRange = $d.Range($changeedPara.Range.Start, $lastPara.Range.End)
$delRange.Delete()

# ------------------------------------------------------------------
# 3. Re-apply paragraph spacing on the now-last paragraph
#    (w:spacing w:before="240" w:after="1440" -> 12pt / 72pt)
# ------------------------------------------------------------------
$targetPara.Format.SpaceBefore = 12
$targetPara.Format.SpaceAfter = 72

# ------------------------------------------------------------------
# 4. Fix the page orientation on the (only) section so it is written
#    out explicitly (w:orient="portrait").
# ------------------------------------------------------------------
$d.Sections(1).PageSetup.Orientation = 0

Write-Output "edit complete"
